# Fruta / hortaliza, semanal
# Insert a new weekly record as row 38 (Fecha serial 44609 = 2022-02-17),
# pushing the existing rows 38-55 down to 39-56.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 38; this shifts rows 38:55 down
# to 39:56 and extends the used range automatically.
$ws.Rows.Item(38).Insert()

# Populate the newly inserted row 38 with the new weekly observation.
$ws.Cells.Item(38, 1).Value = 3
$ws.Cells.Item(38, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(38, 3).Value = "Coquimbo"
$ws.Cells.Item(38, 4).Value = 44609
$ws.Cells.Item(38, 5).Value = 5
$ws.Cells.Item(38, 6).Value = "Fruta"
$ws.Cells.Item(38, 7).Value = 100108
$ws.Cells.Item(38, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(38, 9).Value = 100108004
$ws.Cells.Item(38, 10).Value = "Papaya"
$ws.Cells.Item(38, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(38, 12).Value = "Primera"
$ws.Cells.Item(38, 13).Value = 56
$ws.Cells.Item(38, 14).Value = 23000
$ws.Cells.Item(38, 15).Value = 23000
$ws.Cells.Item(38, 16).Value = 23000
$ws.Cells.Item(38, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(38, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(38, 19).Value = 2300
$ws.Cells.Item(38, 20).Value = 10
